$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -3
$ws.Range("F7").Value = 3
$ws.Range("F9").Value = 2
$ws.Range("F10").Value = -5
$ws.Range("F13").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = -1
